$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Updated fitting parameters
$ws.Range("J2").Value = 0.01409
$ws.Range("K2").Value = 0.1107

# Underline the "-" unit label in K3 (matches J3 styling intent) and move selection there
$ws.Range("K3").Font.Underline = $true
$ws.Range("K3").Select() | Out-Null
